# Update gh-pages to output generated at 456a3b4
# This applies the updated "想去人数" (want-to-go count) figures to the
# "展览" (Exhibitions) sheet and the "全部类型" (All Types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 3,8,9,10,11,12 in column F) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 807
$wsExpo.Range("F8").Value = 4438
$wsExpo.Range("F9").Value = 96
$wsExpo.Range("F10").Value = 4981
$wsExpo.Range("F11").Value = 561
$wsExpo.Range("F12").Value = 1255

# --- Sheet "全部类型" (rows 3,9,10,11,12,13 in column F) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 807
$wsAll.Range("F9").Value = 4438
$wsAll.Range("F10").Value = 96
$wsAll.Range("F11").Value = 4981
$wsAll.Range("F12").Value = 561
$wsAll.Range("F13").Value = 1255
